$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 0.3146067415730337
$ws.Range("D3").Value = 0.6404494382022472
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0.3889655172413793
$ws.Range("I3").Value = 0.1736096388947778
$ws.Range("J3").Value = 0.2134831460674157
$ws.Range("K3").Value = 155.4157303370787

$ws.Range("Q3").Value = 6
$ws.Range("R3").Value = 24
$ws.Range("S3").Value = 65
$ws.Range("T3").Value = 159
$ws.Range("U3").Value = 265
$ws.Range("V3").Value = 630
$ws.Range("W3").Value = 612
$ws.Range("X3").Value = 571
$ws.Range("Y3").Value = 477
$ws.Range("Z3").Value = 371

$ws.Range("AF3").Value = 0.9905659999999999
$ws.Range("AG3").Value = 0.962264
$ws.Range("AH3").Value = 0.897799
$ws.Range("AI3").Value = 0.75
$ws.Range("AJ3").Value = 0.583333
